# Fruta / hortaliza, semanal
# Insert 3 new weekly records at rows 291-293 (pushing the existing
# rows 291-310 down to 294-313) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 291. Doing this three times
# (always at the same row index) shifts the old 291..310 block down to
# 294..313, matching the target dimension A1:T313.
$ws.Rows.Item(291).Insert()
$ws.Rows.Item(291).Insert()
$ws.Rows.Item(291).Insert()

# Common (boilerplate) values shared by every data row in this sheet.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$prodId    = 100101
$producto  = "Berries"
$catId     = 100101001
$categoria = "Arándano (blue)"
$variedad  = "Sin especificar"
$unidad    = "`$/bandeja 2 kilos"
$kgUnidad  = 2

function Set-DataRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $prodId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-DataRow 291 44610 "Primera" 2000 4000 4000 4000 "Provincia de Curicó" 2000
Set-DataRow 292 44610 "Primera" 500  3400 3400 3400 "Región de O'Higgins" 1700
Set-DataRow 293 44610 "Segunda" 250  3200 3200 3200 "Región de O'Higgins" 1600
